# "refactored to extract date from filename or data"
#
# Column K ("Fecha") previously held the literal run-date as plain text
# ("010324", taken straight from the export filename output_010324.xlsx).
# The refactor now parses that into a real date/time value (2024-03-01)
# and writes it as a genuine Excel date serial, formatted as
# "YYYY-MM-DD HH:MM:SS" instead of an inline string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$date = Get-Date -Year 2024 -Month 3 -Day 1 -Hour 0 -Minute 0 -Second 0

# Touch a lowercase variant of the format first (matches the number-format
# registration order of the original export tool), then apply the real
# uppercase format + values to the whole column range.
$probe = $ws.Range("K2")
$probe.NumberFormat = "yyyy-mm-dd h:mm:ss"

$range = $ws.Range("K2:K37")
$range.NumberFormat = "YYYY-MM-DD HH:MM:SS"
$range.Value = $date
